# Insert two new rows at row 253 (shifting old rows 253:363 down to 255:365)
# and populate the two new rows with a new "Perú" / "nueva(o)" price observation
# for the same market/category/variety, dated 2021-09-27 (serial 44466).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("253:254").Insert()

# --- Row 253: "1a nueva(o)" from Perú ---
$ws.Cells.Item(253, 1).Value = 8
$ws.Cells.Item(253, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(253, 3).Value = "Coquimbo"
$ws.Cells.Item(253, 4).Value2 = 44466
$ws.Cells.Item(253, 5).Value = 4
$ws.Cells.Item(253, 6).Value = 100112045
$ws.Cells.Item(253, 7).Value = "Zapallo"
$ws.Cells.Item(253, 8).Value = "Camote"
$ws.Cells.Item(253, 9).Value = "1a nueva(o)"
$ws.Cells.Item(253, 10).Value = 800
$ws.Cells.Item(253, 11).Value = 750
$ws.Cells.Item(253, 12).Value = 800
$ws.Cells.Item(253, 13).Value = 775
$ws.Cells.Item(253, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(253, 15).Value = "Perú"
$ws.Cells.Item(253, 16).Value = 775
$ws.Cells.Item(253, 17).Value = 1
$ws.Cells.Item(253, 18).Value = "Hortaliza"

# --- Row 254: "2a nueva(o)" from Perú ---
$ws.Cells.Item(254, 1).Value = 8
$ws.Cells.Item(254, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(254, 3).Value = "Coquimbo"
$ws.Cells.Item(254, 4).Value2 = 44466
$ws.Cells.Item(254, 5).Value = 4
$ws.Cells.Item(254, 6).Value = 100112045
$ws.Cells.Item(254, 7).Value = "Zapallo"
$ws.Cells.Item(254, 8).Value = "Camote"
$ws.Cells.Item(254, 9).Value = "2a nueva(o)"
$ws.Cells.Item(254, 10).Value = 560
$ws.Cells.Item(254, 11).Value = 650
$ws.Cells.Item(254, 12).Value = 700
$ws.Cells.Item(254, 13).Value = 675
$ws.Cells.Item(254, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(254, 15).Value = "Perú"
$ws.Cells.Item(254, 16).Value = 675
$ws.Cells.Item(254, 17).Value = 1
$ws.Cells.Item(254, 18).Value = "Hortaliza"
